# Daily attendance processing - 2026-01-12 03:40:22
#
# For every data row in the "Recorded By" column (column G), when the
# value is a comma-separated list of recorders whose first entry is
# "System" (case-insensitive), flip the order of the list so "System"
# moves to the end (e.g. "System, foo@bar.com" -> "foo@bar.com, System").
# Rows whose "Recorded By" value does not start with "System" (or has
# only a single entry) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($text -eq $null) {
        continue
    }
    if ($text -eq "") {
        continue
    }
    if ($text.IndexOf(",") -lt 0) {
        continue
    }

    $rawParts = $text.Split(",")
    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    if ($parts[0].ToLower() -ne "system") {
        continue
    }

    $count = $parts.Count
    $reversed = @()
    for ($i = $count - 1; $i -ge 0; $i--) {
        $reversed += $parts[$i]
    }

    $newValue = [string]::Join(", ", $reversed)
    $cell.Value = $newValue
}
